# Scheduled-runner market-data refresh for Pandaemonium_Profits.
# Updates currentAveragePrice* / LevePrice* / LeveProfit* columns (H:N)
# on each job/craft sheet to the latest pulled prices.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1938.1923
$ws.Range("I137").Value = 1513.3611
$ws.Range("J137").Value = 2894.0625
$ws.Range("K137").Value = 4540.0833
$ws.Range("L137").Value = 8682.1875
$ws.Range("M137").Value = -1990.0833
$ws.Range("N137").Value = -13782.1875

$ws.Range("H138").Value = 4154.316
$ws.Range("I138").Value = 2103.45
$ws.Range("J138").Value = 6433.0557
$ws.Range("K138").Value = 6310.349999999999
$ws.Range("L138").Value = 19299.1671
$ws.Range("M138").Value = -1170.349999999999
$ws.Range("N138").Value = -29579.1671


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1132.3572
$ws.Range("I2").Value = 1277.5
$ws.Range("J2").Value = 938.8333
$ws.Range("K2").Value = 1277.5
$ws.Range("L2").Value = 938.8333
$ws.Range("M2").Value = -1164.5
$ws.Range("N2").Value = -1164.8333

$ws.Range("H61").Value = 6236.1387
$ws.Range("I61").Value = 3594.392
$ws.Range("J61").Value = 15859.643
$ws.Range("K61").Value = 3594.392
$ws.Range("L61").Value = 15859.643
$ws.Range("M61").Value = -3382.392
$ws.Range("N61").Value = -16283.643

$ws.Range("H74").Value = 8733.767
$ws.Range("I74").Value = 5715.3213
$ws.Range("J74").Value = 50992
$ws.Range("K74").Value = 5715.3213
$ws.Range("L74").Value = 50992
$ws.Range("M74").Value = -4841.3213
$ws.Range("N74").Value = -52740

$ws.Range("H77").Value = 8733.767
$ws.Range("I77").Value = 5715.3213
$ws.Range("J77").Value = 50992
$ws.Range("K77").Value = 28576.6065
$ws.Range("L77").Value = 254960
$ws.Range("M77").Value = -24208.6065
$ws.Range("N77").Value = -263696

$ws.Range("H116").Value = 1132.3572
$ws.Range("I116").Value = 1277.5
$ws.Range("J116").Value = 938.8333
$ws.Range("K116").Value = 1277.5
$ws.Range("L116").Value = 938.8333
$ws.Range("M116").Value = 1016.5
$ws.Range("N116").Value = -5526.8333

$ws.Range("H132").Value = 5381.857
$ws.Range("I132").Value = 2020.5714
$ws.Range("J132").Value = 7622.7144
$ws.Range("K132").Value = 6061.7142
$ws.Range("L132").Value = 22868.1432
$ws.Range("M132").Value = -3531.7142
$ws.Range("N132").Value = -27928.1432

$ws.Range("H136").Value = 6236.1387
$ws.Range("I136").Value = 3594.392
$ws.Range("J136").Value = 15859.643
$ws.Range("K136").Value = 10783.176
$ws.Range("L136").Value = 47578.929
$ws.Range("M136").Value = -8233.175999999999
$ws.Range("N136").Value = -52678.929


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1132.3572
$ws.Range("I3").Value = 1277.5
$ws.Range("J3").Value = 938.8333
$ws.Range("K3").Value = 1277.5
$ws.Range("L3").Value = 938.8333
$ws.Range("M3").Value = -1163.5
$ws.Range("N3").Value = -1166.8333

$ws.Range("H134").Value = 17211.785
$ws.Range("I134").Value = 1672.5625
$ws.Range("J134").Value = 61087.234
$ws.Range("K134").Value = 5017.6875
$ws.Range("L134").Value = 183261.702
$ws.Range("M134").Value = -2482.6875
$ws.Range("N134").Value = -188331.702


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2119.2104
$ws.Range("I31").Value = 1550.7906
$ws.Range("J31").Value = 3865.0715
$ws.Range("K31").Value = 1550.7906
$ws.Range("L31").Value = 3865.0715
$ws.Range("M31").Value = -1255.7906
$ws.Range("N31").Value = -4455.0715

$ws.Range("H34").Value = 2119.2104
$ws.Range("I34").Value = 1550.7906
$ws.Range("J34").Value = 3865.0715
$ws.Range("K34").Value = 1550.7906
$ws.Range("L34").Value = 3865.0715
$ws.Range("M34").Value = -1348.7906
$ws.Range("N34").Value = -4269.0715

$ws.Range("H69").Value = 0
$ws.Range("I69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("M69").ClearContents()

$ws.Range("H70").Value = 71000
$ws.Range("J70").Value = 71000
$ws.Range("L70").Value = 71000
$ws.Range("N70").Value = -71630

$ws.Range("H72").Value = 0
$ws.Range("I72").Value = 0
$ws.Range("K72").Value = 0
$ws.Range("M72").ClearContents()

$ws.Range("H73").Value = 71000
$ws.Range("J73").Value = 71000
$ws.Range("L73").Value = 71000
$ws.Range("N73").Value = -73184

$ws.Range("H75").Value = 0
$ws.Range("I75").Value = 0
$ws.Range("K75").Value = 0
$ws.Range("M75").ClearContents()

$ws.Range("H78").Value = 0
$ws.Range("I78").Value = 0
$ws.Range("K78").Value = 0
$ws.Range("M78").ClearContents()

$ws.Range("H81").Value = 40328
$ws.Range("J81").Value = 40328
$ws.Range("L81").Value = 40328
$ws.Range("N81").Value = -42324

$ws.Range("H84").Value = 40328
$ws.Range("J84").Value = 40328
$ws.Range("L84").Value = 120984
$ws.Range("N84").Value = -130968

$ws.Range("H132").Value = 3757.5
$ws.Range("I132").Value = 4276.5264
$ws.Range("J132").Value = 2861
$ws.Range("K132").Value = 12829.5792
$ws.Range("L132").Value = 8583
$ws.Range("M132").Value = -10299.5792
$ws.Range("N132").Value = -13643


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 6945252
$ws.Range("I5").Value = 463.7143
$ws.Range("J5").Value = 16667955
$ws.Range("K5").Value = 1391.1429
$ws.Range("L5").Value = 50003865
$ws.Range("M5").Value = -1279.1429
$ws.Range("N5").Value = -50004089

$ws.Range("H135").Value = 6945252
$ws.Range("I135").Value = 463.7143
$ws.Range("J135").Value = 16667955
$ws.Range("K135").Value = 4173.428699999999
$ws.Range("L135").Value = 150011595
$ws.Range("M135").Value = -1638.428699999999
$ws.Range("N135").Value = -150016665


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 7776.1
$ws.Range("I122").Value = 18333.666
$ws.Range("J122").Value = 3251.4285
$ws.Range("K122").Value = 55000.99800000001
$ws.Range("L122").Value = 9754.2855
$ws.Range("M122").Value = -52550.99800000001
$ws.Range("N122").Value = -14654.2855

$ws.Range("H132").Value = 4065.3333
$ws.Range("I132").Value = 1711.0256
$ws.Range("J132").Value = 34671.332
$ws.Range("K132").Value = 5133.0768
$ws.Range("L132").Value = 104013.996
$ws.Range("M132").Value = -2603.0768
$ws.Range("N132").Value = -109073.996


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3378.4395
$ws.Range("I132").Value = 3445.0588
$ws.Range("J132").Value = 3151.9333
$ws.Range("K132").Value = 10335.1764
$ws.Range("L132").Value = 9455.7999
$ws.Range("M132").Value = -7805.1764
$ws.Range("N132").Value = -14515.7999


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1357.8235
$ws.Range("I132").Value = 695.129
$ws.Range("J132").Value = 2385
$ws.Range("K132").Value = 2085.387
$ws.Range("L132").Value = 7155
$ws.Range("M132").Value = 444.6129999999998
$ws.Range("N132").Value = -12215
